$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Partnersysteme" group header (row 23) to make
# room for the new "Anzahl Partnersysteme" requirement (id 3100). This shifts
# every row from 23 onward down by one.
$ws.Rows("23:23").Insert()

# Fill in the new requirement row. Column order matches how the values were
# originally authored (Titel, then Beschreibung, then Quelle) so that new
# shared-string entries come out in the same order as the source workbook.
$ws.Range("B23").Value = 3100
$ws.Range("C23").Value = "B"
$ws.Range("D23").Value = "f"
$ws.Range("F23").Value = "Anzahl Partnersysteme"
$ws.Range("G23").Value = "Die Anzahl an zusätzlichen Partnersystemen, die für die Lösung benötigt werden, soll so gering wie möglich gehalten werden."
$ws.Range("E23").Value = "K"

$ws.Rows("23:23").RowHeight = 30

# Restore the scroll/selection state recorded for the edited sheet.
$ws.Range("F27").Select()
